$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Row 67: the pending XRP/USDT Sell order (txn code "XRP/USDT0000009") is
# now finalized ("DONE"), gets a finalized date, a fee, a profit % and a
# transaction duration.
# ---------------------------------------------------------------------------
$ws.Range("H67").Value2 = "DONE"
$ws.Range("I67").Value2 = 42878.568287037036

$ws.Range("J67").Value2 = "0.06823749 USDT (0.15%)"

$ws.Range("K67").Value2 = "     ~0.2%"
$k67 = $ws.Range("K67").Characters(6, 5)
$k67.Font.Color = 5287936

$ws.Range("L67").Value2 = "2 day"

# The transaction value text for row 67 also got corrected.
$ws.Range("E67").Value2 = "         0.340  USDT"

# ---------------------------------------------------------------------------
# Row 68: new Buy order (XRP/USDT0000010), already finalized.
# ---------------------------------------------------------------------------
$ws.Range("A68").Value2 = "             2017-05-24 10:47:02"
$ws.Range("I66").Copy() | Out-Null
$ws.Range("A68").PasteSpecial(-4122) | Out-Null
$ws.Range("A68").Value2 = "             2017-05-24 10:47:02"

$ws.Range("B68").Value2 = "            Buy"
$b68 = $ws.Range("B68").Characters(13, 3)
$b68.Font.Color = 5287936

$ws.Range("C68").Value2 = "        XRP"

$ws.Range("D67").Copy() | Out-Null
$ws.Range("D68").PasteSpecial(-4122) | Out-Null
$ws.Range("D68").Value2 = 0.3155

$ws.Range("E68").Value2 = "         0.316  USDT"
$ws.Range("F68").Value2 = "         144 XRP"
$ws.Range("G68").Value2 = " XRP/USDT0000010"

$ws.Range("H67").Copy() | Out-Null
$ws.Range("H68").PasteSpecial(-4122) | Out-Null
$ws.Range("H68").Value2 = "DONE"

$ws.Range("I66").Copy() | Out-Null
$ws.Range("I68").PasteSpecial(-4122) | Out-Null
$ws.Range("I68").Value2 = 42879.449328703704

$ws.Range("J68").Value2 = "0.21674578 XRP (0.15%)"
$ws.Range("K68").Value2 = "     "

$ws.Rows.Item(68).RowHeight = 14.25

# ---------------------------------------------------------------------------
# Row 69: new Sell order (same txn code, "XRP/USDT0000010"), still pending.
# ---------------------------------------------------------------------------
$ws.Range("I66").Copy() | Out-Null
$ws.Range("A69").PasteSpecial(-4122) | Out-Null
$ws.Range("A69").Value2 = "                2017-05-24 10:47:02"

$ws.Range("B69").Value2 = "            Sell"
$b69 = $ws.Range("B69").Characters(13, 4)
$b69.Font.Color = 255

$ws.Range("C69").Value2 = "        XRP"

$ws.Range("D67").Copy() | Out-Null
$ws.Range("D69").PasteSpecial(-4122) | Out-Null
$ws.Range("D69").Value2 = 0.319

$ws.Range("E69").Value2 = "         0.342  USDT"
$ws.Range("F69").Value2 = "         144 XRP"
$ws.Range("G69").Value2 = " XRP/USDT0000010"

$ws.Range("H67").Copy() | Out-Null
$ws.Range("H69").PasteSpecial(-4122) | Out-Null
$ws.Range("H69").Value2 = "IN PROGRESS"

$ws.Range("I66").Copy() | Out-Null
$ws.Range("I69").PasteSpecial(-4122) | Out-Null

$ws.Range("K69").Value2 = "    "

$ws.Rows.Item(69).RowHeight = 14.25

# ---------------------------------------------------------------------------
# View state: the author scrolled down and landed the selection on D69.
# ---------------------------------------------------------------------------
$ws.Range("D69").Select()
